$d = $word.ActiveDocument

# Locate the paragraph ending in "Thomson Pioneira (2008)." and the paragraph
# containing the "Creative Commons Attribution" copyright text, then remove
# everything between the end of the former and the end of the latter
# (i.e. the blank paragraph, the "Ver no Jupiter..." paragraph, and the
# copyright paragraph), leaving the "Thomson Pioneira (2008)." paragraph and
# the blank paragraph that follows the deleted block intact.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Thomson Pioneira*") {
        $startPara = $p
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.End, $endPara.Range.End)
    $r.Delete()
}
